$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# IS29 Inkjet ink rows (4-7): swap Black<->Magenta and Cyan<->Yellow labels
$ws.Range("D4").Value = "Magenta - IS29 Inkjet - "
$ws.Range("D5").Value = "Black - IS29 Inkjet - "
$ws.Range("D6").Value = "Yellow - IS29 Inkjet - "
$ws.Range("D7").Value = "Cyan - IS29 Inkjet - "

# Digital ink rows (9-12): rotate color labels
$ws.Range("D9").Value = "Cyan - Digital - "
$ws.Range("D10").Value = "Black - Digital - "
$ws.Range("D11").Value = "Magenta - Digital - "
$ws.Range("D12").Value = "Yellow - Digital - "
